# The Title, Author and Abstract paragraphs each had their text split
# across many runs (one run per word, plus separate single-space runs
# in between). The edit collapses each of those paragraphs down to a
# single run containing the full (unchanged) text. Using Find/Replace
# with the exact same text as both the search and replacement target
# causes Word to rewrite the matched range as one consolidated run,
# which is exactly the structural change the diff describes.
$d = $word.ActiveDocument

$d.Content.Find.Execute(
    "Questions: Introduction to hypothesis testing", $true, $false, $false, $false, $false,
    $true, 1, $false, "Questions: Introduction to hypothesis testing", 2)

$d.Content.Find.Execute(
    "Ellie Trace", $true, $false, $false, $false, $false,
    $true, 1, $false, "Ellie Trace", 2)

$d.Content.Find.Execute(
    "A selection of questions for the study guide on introduction to hypothesis testing.", $true, $false, $false, $false, $false,
    $true, 1, $false, "A selection of questions for the study guide on introduction to hypothesis testing.", 2)
